# Updated cryptos list price (D) and volume(1h) (E) values for rows 2-51.
# Values in column D that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# original inlineStr cell contents) instead of re-parsing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.573.73"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.983.53"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'587.41"
$ws.Range("E5").Value = "  +12.28%  "
$ws.Range("D6").Value = "'152.09"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.748"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'53.44"
$ws.Range("E11").Value = "  +6.27%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'10.81"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "4.623.23"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "3.989.62"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  +9.33%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'20.43"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'0.132"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "72.487.42"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'429.50"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'4.72"
$ws.Range("E22").Value = "  +14.51%  "
$ws.Range("D23").Value = "'95.67"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'4.50"
$ws.Range("E25").Value = "  +22.24%  "
$ws.Range("D26").Value = "'14.24"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "'11.25"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").Value = "'10.59"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = "  +6.35%  "
$ws.Range("D32").Value = "'50.04"
$ws.Range("E32").Value = "  +3.66%  "
$ws.Range("D33").Value = "'13.48"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "'0.131"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "'682.27"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'68.98"
$ws.Range("E36").Value = "  +5.99%  "
$ws.Range("D37").Value = "'0.437"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "0.0₃0853"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'11.10"
$ws.Range("E42").Value = "  +12.77%  "
$ws.Range("D43").Value = "'3.30"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("D47").Value = "'0.148"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'3.36"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'3.46"
$ws.Range("E49").Value = "  +7.18%  "
$ws.Range("D50").Value = "'2.99"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "'2.15"
$ws.Range("E51").Value = "  +7.53%  "
